$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 635
$ws1.Range("F3").Value = 2218
$ws1.Range("F4").Value = 92
$ws1.Range("F5").Value = 13286
$ws1.Range("F11").Value = 995
$ws1.Range("F12").Value = 13808
$ws1.Range("F13").Value = 14446
$ws1.Range("F15").Value = 172
$ws1.Range("F22").Value = 1101
$ws1.Range("F24").Value = 58
$ws1.Range("F25").Value = 5493
$ws1.Range("F27").Value = 642
$ws1.Range("F28").Value = 347
$ws1.Range("F30").Value = 98

# Sheet "全部类型" (sheet4) updates to column F (想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 635
$ws4.Range("F3").Value = 2218
$ws4.Range("F4").Value = 92
$ws4.Range("F5").Value = 13286
$ws4.Range("F12").Value = 995
$ws4.Range("F13").Value = 13808
$ws4.Range("F14").Value = 14446
$ws4.Range("F16").Value = 172
$ws4.Range("F23").Value = 1101
$ws4.Range("F25").Value = 58
$ws4.Range("F26").Value = 5493
$ws4.Range("F28").Value = 642
$ws4.Range("F29").Value = 347
$ws4.Range("F31").Value = 98
